$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E (shifts old E:Q -> F:R, carrying values/formats along)
$ws.Columns("E:E").Insert()

# The insert left "ghost" formatted-but-empty cells at E3/E4 (copied from the old E
# column's date-input style). The target layout has no cell there at all, so fully
# clear them (contents + formatting), leaving the new F3/F4 (shifted from old E3/E4)
# as the only populated cells in that area.
$ws.Range("E3").Clear()
$ws.Range("E4").Clear()

# New header cell for the inserted column - copy formatting from the header row's
# existing style (row 10 header cells all share one style) and set the caption.
$ws.Range("D10").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("E10").Value = "Billing Responsible"

# The old "Billing Responsible" header (shifted from P10 to Q10 by the column
# insert) is being renamed to "Invoice Responsible".
$ws.Range("Q10").Value = "Invoice Responsible"

# Column formatting: the newly inserted column E needs an explicit width/style
# (a plain text/general column, matching style used by most data columns).
$ws.Columns("E:E").ColumnWidth = 26.82
$ws.Range("E1").NumberFormat = $ws.Range("F1").NumberFormat

# Column Q (shifted from old column P) changes from the general style to the
# date-format style, with a slightly adjusted width.
$ws.Columns("Q:Q").ColumnWidth = 26.95
$ws.Range("Q1").NumberFormat = "DD/MM/YYYY"
